$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full Lamc2 -> Itga3 edge table across the ECs / FAPs / sCs cluster cross-product
$rows = @(
    @{ A="ECs"; D="ECs"; E=2; F=0.6666666666666666; G=0.4421816666666667; H=1.326545; I=0.06026482003168283; J=0.06026482003168283; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=3.656474782974444; R=32.90827304677; S=0.04047868959430988; T=0.04047868959430988 },
    @{ A="ECs"; D="FAPs"; E=2; F=0.6666666666666666; G=0.4421816666666667; H=1.326545; I=0.06026482003168283; J=0.06026482003168283; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.05363398307666667; R=0.48270584769; S=0.0005937503966322398; T=0.0005937503966322398 },
    @{ A="ECs"; D="sCs"; E=2; F=0.6666666666666666; G=0.4421816666666667; H=1.326545; I=0.06026482003168283; J=0.06026482003168283; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=1.733664166196111; R=15.602977495765; S=0.01919238004074071; T=0.01919238004074071 },
    @{ A="FAPs"; D="ECs"; E=3; F=1; G=6.369908666666666; H=19.109726; I=0.8681531333236113; J=0.8681531333236113; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=52.67384915592844; R=474.0646424033559; S=0.5831213166430937; T=0.5831213166430937 },
    @{ A="FAPs"; D="FAPs"; E=3; F=1; G=6.369908666666666; H=19.109726; I=0.8681531333236113; J=0.8681531333236113; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.7726317018146667; R=6.953685316331999; S=0.00855335280147558; T=0.00855335280147558 },
    @{ A="FAPs"; D="sCs"; E=3; F=1; G=6.369908666666666; H=19.109726; I=0.8681531333236113; J=0.8681531333236113; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=24.97453700554911; R=224.770833049942; S=0.276478463879042; T=0.276478463879042 },
    @{ A="sCs"; D="ECs"; E=3; F=1; G=0.5252196666666666; H=1.575659; I=0.07158204664470585; J=0.07158204664470584; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=4.34313001071711; R=39.08817009645399; S=0.04808024723434237; T=0.04808024723434236 },
    @{ A="sCs"; D="FAPs"; E=3; F=1; G=0.5252196666666666; H=1.575659; I=0.07158204664470585; J=0.07158204664470584; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.06370599424866666; R=0.5733539482379999; S=0.000705251730025863; T=0.0007052517300258628 },
    @{ A="sCs"; D="sCs"; E=3; F=1; G=0.5252196666666666; H=1.575659; I=0.07158204664470585; J=0.07158204664470584; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=2.059231723344777; R=18.533085510103; S=0.02279654768033762; T=0.02279654768033761 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row.A      # Sending cluster
    $ws.Cells.Item($r, 2).Value  = "Lamc2"     # Ligand symbol
    $ws.Cells.Item($r, 3).Value  = "Itga3"     # Receptor symbol
    $ws.Cells.Item($r, 4).Value  = $row.D      # Target cluster
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}
